$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
  "2021-10-05 13:38:49.332505",
  "2021-10-05 13:38:49.332514",
  "2021-10-05 13:38:49.332518",
  "2021-10-05 13:38:49.332520",
  "2021-10-05 13:38:49.332523",
  "2021-10-05 13:38:49.332525",
  "2021-10-05 13:38:49.332527",
  "2021-10-05 13:38:49.332529",
  "2021-10-05 13:38:49.332532",
  "2021-10-05 13:38:49.332534",
  "2021-10-05 13:38:49.332537",
  "2021-10-05 13:38:49.332539",
  "2021-10-05 13:38:49.332541",
  "2021-10-05 13:38:49.332543",
  "2021-10-05 13:38:49.332546",
  "2021-10-05 13:38:49.332548",
  "2021-10-05 13:38:49.332550",
  "2021-10-05 13:38:49.332553",
  "2021-10-05 13:38:49.332555",
  "2021-10-05 13:38:49.332557",
  "2021-10-05 13:38:49.332560",
  "2021-10-05 13:38:49.332562",
  "2021-10-05 13:38:49.332564",
  "2021-10-05 13:38:49.332566",
  "2021-10-05 13:38:49.332569",
  "2021-10-05 13:38:49.332572",
  "2021-10-05 13:38:49.332574",
  "2021-10-05 13:38:49.332576",
  "2021-10-05 13:38:49.332578",
  "2021-10-05 13:38:49.332581",
  "2021-10-05 13:38:49.332583",
  "2021-10-05 13:38:49.332585",
  "2021-10-05 13:38:49.332588",
  "2021-10-05 13:38:49.332590",
  "2021-10-05 13:38:49.332593",
  "2021-10-05 13:38:49.332595",
  "2021-10-05 13:38:49.332597",
  "2021-10-05 13:38:49.332599",
  "2021-10-05 13:38:49.332602",
  "2021-10-05 13:38:49.332604",
  "2021-10-05 13:38:49.332606",
  "2021-10-05 13:38:49.332609",
  "2021-10-05 13:38:49.332611",
  "2021-10-05 13:38:49.332613",
  "2021-10-05 13:38:49.332615",
  "2021-10-05 13:38:49.332618",
  "2021-10-05 13:38:49.332620",
  "2021-10-05 13:38:49.332622",
  "2021-10-05 13:38:49.332624",
  "2021-10-05 13:38:49.332627",
  "2021-10-05 13:38:49.332629",
  "2021-10-05 13:38:49.332631",
  "2021-10-05 13:38:49.332634",
  "2021-10-05 13:38:49.332636",
  "2021-10-05 13:38:49.332638",
  "2021-10-05 13:38:49.332641",
  "2021-10-05 13:38:49.332643",
  "2021-10-05 13:38:49.332645",
  "2021-10-05 13:38:49.332647",
  "2021-10-05 13:38:49.332650",
  "2021-10-05 13:38:49.332652",
  "2021-10-05 13:38:49.332654",
  "2021-10-05 13:38:49.332656",
  "2021-10-05 13:38:49.332659",
  "2021-10-05 13:38:49.332662",
  "2021-10-05 13:38:49.332664",
  "2021-10-05 13:38:49.332667",
  "2021-10-05 13:38:49.332669",
  "2021-10-05 13:38:49.332671",
  "2021-10-05 13:38:49.332673",
  "2021-10-05 13:38:49.332676",
  "2021-10-05 13:38:49.332678",
  "2021-10-05 13:38:49.332680",
  "2021-10-05 13:38:49.332683",
  "2021-10-05 13:38:49.332685",
  "2021-10-05 13:38:49.332687",
  "2021-10-05 13:38:49.332691",
  "2021-10-05 13:38:49.332694",
  "2021-10-05 13:38:49.332697",
  "2021-10-05 13:38:49.332699",
  "2021-10-05 13:38:49.332701",
  "2021-10-05 13:38:49.332704",
  "2021-10-05 13:38:49.332706",
  "2021-10-05 13:38:49.332708",
  "2021-10-05 13:38:49.332710",
  "2021-10-05 13:38:49.332713",
  "2021-10-05 13:38:49.332715",
  "2021-10-05 13:38:49.332717",
  "2021-10-05 13:38:49.332719",
  "2021-10-05 13:38:49.332722",
  "2021-10-05 13:38:49.332724",
  "2021-10-05 13:38:49.332726",
  "2021-10-05 13:38:49.332730",
  "2021-10-05 13:38:49.332732",
  "2021-10-05 13:38:49.332734",
  "2021-10-05 13:38:49.332737",
  "2021-10-05 13:38:49.332739",
  "2021-10-05 13:38:49.332741",
  "2021-10-05 13:38:49.332743",
  "2021-10-05 13:38:49.332746",
  "2021-10-05 13:38:49.332748",
  "2021-10-05 13:38:49.332750",
  "2021-10-05 13:38:49.332752",
  "2021-10-05 13:38:49.332755",
  "2021-10-05 13:38:49.332757",
  "2021-10-05 13:38:49.332759",
  "2021-10-05 13:38:49.332761",
  "2021-10-05 13:38:49.332764",
  "2021-10-05 13:38:49.332768",
  "2021-10-05 13:38:49.332771",
  "2021-10-05 13:38:49.332773",
  "2021-10-05 13:38:49.332775",
  "2021-10-05 13:38:49.332778",
  "2021-10-05 13:38:49.332780"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
